$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.427.22'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '1.800.71'
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.52'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.600'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '36.31'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.88%  '
$ws.Range("E9").Value = '  -2.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0677'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.67%  '
$ws.Range("E11").Value = '  +1.24%  '
$ws.Range("D12").Value = '2.057.55'
$ws.Range("E12").Value = '  -1.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.25'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("D14").Value = '1.790.38'
$ws.Range("E14").Value = '  -1.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.632'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.28%  '
$ws.Range("D16").Value = '34.381.02'
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.43'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.73'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.91'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").Value = '0.0₃0773'
$ws.Range("E20").Value = '  -3.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.35'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.08'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.70'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.90'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.35'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.121'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.84%  '
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("E30").Value = '  -1.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.79'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.90'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0515'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.80%  '
$ws.Range("E34").Value = '  -3.82%  '
$ws.Range("D35").Value = '1.361.71'
$ws.Range("E35").Value = '  -2.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.648'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.45%  '
$ws.Range("E37").Value = '  -1.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.36'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -7.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0187'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.41'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.79'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '81.10'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.937'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("E44").Value = '  +4.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.24'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0497'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.20%  '
$ws.Range("D47").Value = '1.961.15'
$ws.Range("E47").Value = '  -0.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.77'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.66%  '
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.92'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.32%  '
$ws.Range("E51").Value = '  -7.41%  '
